$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "R1-R2, R6-R8, R10-R13, R16"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "1 x 10 = 10 pt(s)"

# Row 3
$ws.Range("A3").Value = "R9, R17-R20"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "2 x 5 = 10 pt(s)"

# Row 4
$ws.Range("A4").Value = "R3, R15, R22"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "3 x 3 = 9 pt(s)"

# Row 5
$ws.Range("A5").Value = "R4-R5, R14, R21, R23"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "4 x 5 = 20 pt(s)"

# Row 6
$ws.Range("A6").Value = "S1-S7"
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = "3 x 7 = 21 pt(s)"
